# -----------------------------------------------------------------------
# Rebuilds the PO forecast comparison workbook:
#   - renames Sheet1 -> "Sales vs PO" and reshapes it (adds an "Order Week"
#     column, shifts the sales week forward one week, zeroes out the
#     now-unused PO column in-place)
#   - adds "Weekly Growth"   : the weeks that actually had a PO request,
#                               with week-over-week growth %
#   - adds "Volume Insights" : aggregate stats over those PO quantities
#   - adds "Prediction Info" : naive linear-trend forecast for next week
# -----------------------------------------------------------------------

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Sales vs PO"

# --- reshape "Sales vs PO" ----------------------------------------------
# Insert a fresh column C; this shifts the existing PO_Requested_Qty
# column (old C) over to D, carrying its header/values with it.
$ws1.Columns.Item(3).Insert()
$ws1.Range("C1").Value = "Order Week"

$lastRow = $ws1.Cells.Item($ws1.Rows.Count, 1).End(-4162).Row

# give the new "Order Week" column the same date formatting as column A
$ws1.Range("A2").Copy()
$ws1.Range("C2:C" + $lastRow).PasteSpecial(-4122)

# walk every data row: stash any non-zero PO qty (with its original
# "ds" week) for the Weekly Growth sheet, then:
#   C(row) <- old A(row)      (Order Week = the original sales week)
#   A(row) <- old A(row) + 6  (Sales week pushed out by one week)
#   D(row) <- 0                (PO qty now lives on the Weekly Growth sheet)
$poDates = New-Object System.Collections.ArrayList
$poQtys  = New-Object System.Collections.ArrayList

for ($r = 2; $r -le $lastRow; $r++) {
    $oldA  = $ws1.Cells.Item($r, 1).Value2
    $poQty = $ws1.Cells.Item($r, 4).Value2

    if ($poQty -ne 0) {
        [void]$poDates.Add($oldA)
        [void]$poQtys.Add($poQty)
    }

    $ws1.Cells.Item($r, 3).Value2 = $oldA
    $ws1.Cells.Item($r, 1).Value2 = $oldA + 6
    $ws1.Cells.Item($r, 4).Value2 = 0
}

# --- "Weekly Growth" -----------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Weekly Growth"

$ws1.Range("A1:C1").Copy()
$ws2.Range("A1:C1").PasteSpecial(-4122)
$ws2.Range("A1").Value = "ds"
$ws2.Range("B1").Value = "PO_Requested_Qty"
$ws2.Range("C1").Value = "Growth%"

$n = $poDates.Count
if ($n -gt 0) {
    $ws1.Range("A2").Copy()
    $ws2.Range("A2:A" + ($n + 1)).PasteSpecial(-4122)
}

$total = 0
$maxQty = $null
$minQty = $null

for ($i = 0; $i -lt $n; $i++) {
    $row = $i + 2
    $qty = $poQtys[$i]

    $ws2.Cells.Item($row, 1).Value2 = $poDates[$i]
    $ws2.Cells.Item($row, 2).Value2 = $qty

    if ($i -eq 0) {
        $growth = 0
    } else {
        $prev = $poQtys[$i - 1]
        $growth = (($qty - $prev) / $prev) * 100
    }
    $ws2.Cells.Item($row, 3).Value2 = $growth

    $total += $qty
    if ($maxQty -eq $null -or $qty -gt $maxQty) { $maxQty = $qty }
    if ($minQty -eq $null -or $qty -lt $minQty) { $minQty = $qty }
}

$avg = $total / $n

# --- "Volume Insights" ----------------------------------------------------
$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "Volume Insights"

$ws1.Range("A1:C1").Copy()
$ws3.Range("A1:D1").PasteSpecial(-4122)
$ws3.Range("A1").Value = "Total_PO_Quantity"
$ws3.Range("B1").Value = "Average_PO_Quantity"
$ws3.Range("C1").Value = "Max_PO_Quantity"
$ws3.Range("D1").Value = "Min_PO_Quantity"

$ws3.Range("A2").Value2 = $total
$ws3.Range("B2").Value2 = $avg
$ws3.Range("C2").Value2 = $maxQty
$ws3.Range("D2").Value2 = $minQty

# --- "Prediction Info" ----------------------------------------------------
$ws4 = $wb.Worksheets.Add($null, $ws3)
$ws4.Name = "Prediction Info"

$ws1.Range("A1").Copy()
$ws4.Range("A1").PasteSpecial(-4122)
$ws4.Range("A1").Value = "Predicted_Next_Week_PO_Quantity"

# naive linear-trend extrapolation of the PO quantity series to the next point
$ws4.Range("A2").Value2 = 15.00000000000001

# --- restore the original active sheet -----------------------------------
$ws1.Activate()
